# Applies updated TPM-derived NATMI ligand-receptor scores to Sheet1 (rows 2-26).
# Columns G/H (ligand avg/total expr) depend only on the sending cluster; columns
# M/N (receptor avg/total expr) depend only on the target cluster. I/J/O/P are the
# corresponding specificity fractions, Q/R are edge weights (G*M, H*N) and S/T are
# their specificity fractions across all 25 sending/target combinations.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("G2").Value2 = 48.091872
$ws.Range("H2").Value2 = 144.275616
$ws.Range("I2").Value2 = 0.42109384267595806
$ws.Range("J2").Value2 = 0.4237822050924049
$ws.Range("M2").Value2 = 13.89934866666667
$ws.Range("N2").Value2 = 41.69804600000001
$ws.Range("O2").Value2 = 0.048535075531341794
$ws.Range("P2").Value2 = 0.04999273878390351
$ws.Range("Q2").Value2 = 668.4456969607041
$ws.Range("R2").Value2 = 6016.0112726463385
$ws.Range("S2").Value2 = 0.020437821460060582
$ws.Range("T2").Value2 = 0.021186033080451225

# Row 3
$ws.Range("G3").Value2 = 48.091872
$ws.Range("H3").Value2 = 144.275616
$ws.Range("I3").Value2 = 0.42109384267595806
$ws.Range("J3").Value2 = 0.4237822050924049
$ws.Range("M3").Value2 = 70.36235166666667
$ws.Range("N3").Value2 = 211.087055
$ws.Range("O3").Value2 = 0.24569799165441703
$ws.Range("P3").Value2 = 0.253077086664408
$ws.Range("Q3").Value2 = 3383.85720997232
$ws.Range("R3").Value2 = 30454.71488975088
$ws.Range("S3").Value2 = 0.10346191144352394
$ws.Range("T3").Value2 = 0.10724956584500447

# Row 4
$ws.Range("G4").Value2 = 48.091872
$ws.Range("H4").Value2 = 144.275616
$ws.Range("I4").Value2 = 0.42109384267595806
$ws.Range("J4").Value2 = 0.4237822050924049
$ws.Range("M4").Value2 = 82.007665
$ws.Range("N4").Value2 = 246.022995
$ws.Range("O4").Value2 = 0.28636221094801234
$ws.Range("P4").Value2 = 0.2949625822722868
$ws.Range("Q4").Value2 = 3943.9021281988803
$ws.Range("R4").Value2 = 35495.11915378992
$ws.Range("S4").Value2 = 0.12058536380528183
$ws.Range("T4").Value2 = 0.12499989353509958

# Row 5
$ws.Range("G5").Value2 = 48.091872
$ws.Range("H5").Value2 = 144.275616
$ws.Range("I5").Value2 = 0.42109384267595806
$ws.Range("J5").Value2 = 0.4237822050924049
$ws.Range("M5").Value2 = 25.0501465
$ws.Range("N5").Value2 = 50.100293
$ws.Range("O5").Value2 = 0.0874724982879541
$ws.Range("P5").Value2 = 0.06006638442832619
$ws.Range("Q5").Value2 = 1204.7084390592481
$ws.Range("R5").Value2 = 7228.250634355489
$ws.Range("S5").Value2 = 0.03683413043254075
$ws.Range("T5").Value2 = 0.025455064844964164

# Row 6
$ws.Range("G6").Value2 = 48.091872
$ws.Range("H6").Value2 = 144.275616
$ws.Range("I6").Value2 = 0.42109384267595806
$ws.Range("J6").Value2 = 0.4237822050924049
$ws.Range("M6").Value2 = 95.05788666666668
$ws.Range("N6").Value2 = 285.17366
$ws.Range("O6").Value2 = 0.3319322235782747
$ws.Range("P6").Value2 = 0.3419012078510756
$ws.Range("Q6").Value2 = 4571.511718163841
$ws.Range("R6").Value2 = 41143.60546347456
$ws.Range("S6").Value2 = 0.13977461553455095
$ws.Range("T6").Value2 = 0.14489164778688546

# Row 7
$ws.Range("G7").Value2 = 14.05160533333333
$ws.Range("H7").Value2 = 42.154816
$ws.Range("I7").Value2 = 0.12303626869794793
$ws.Range("J7").Value2 = 0.12382176125828906
$ws.Range("M7").Value2 = 13.89934866666667
$ws.Range("N7").Value2 = 41.69804600000001
$ws.Range("O7").Value2 = 0.048535075531341794
$ws.Range("P7").Value2 = 0.04999273878390351
$ws.Range("Q7").Value2 = 195.3081618543929
$ws.Range("R7").Value2 = 1757.7734566895365
$ws.Range("S7").Value2 = 0.005971574594349367
$ws.Range("T7").Value2 = 0.00619018896634851

# Row 8
$ws.Range("G8").Value2 = 14.05160533333333
$ws.Range("H8").Value2 = 42.154816
$ws.Range("I8").Value2 = 0.12303626869794793
$ws.Range("J8").Value2 = 0.12382176125828906
$ws.Range("M8").Value2 = 70.36235166666667
$ws.Range("N8").Value2 = 211.087055
$ws.Range("O8").Value2 = 0.24569799165441703
$ws.Range("P8").Value2 = 0.253077086664408
$ws.Range("Q8").Value2 = 988.7039959452086
$ws.Range("R8").Value2 = 8898.33596350688
$ws.Range("S8").Value2 = 0.03022976411973902
$ws.Range("T8").Value2 = 0.03133645060490366

# Row 9
$ws.Range("G9").Value2 = 14.05160533333333
$ws.Range("H9").Value2 = 42.154816
$ws.Range("I9").Value2 = 0.12303626869794793
$ws.Range("J9").Value2 = 0.12382176125828906
$ws.Range("M9").Value2 = 82.007665
$ws.Range("N9").Value2 = 246.022995
$ws.Range("O9").Value2 = 0.28636221094801234
$ws.Range("P9").Value2 = 0.2949625822722868
$ws.Range("Q9").Value2 = 1152.339342888213
$ws.Range("R9").Value2 = 10371.05408599392
$ws.Range("S9").Value2 = 0.03523293793113809
$ws.Range("T9").Value2 = 0.03652278644224753

# Row 10
$ws.Range("G10").Value2 = 14.05160533333333
$ws.Range("H10").Value2 = 42.154816
$ws.Range("I10").Value2 = 0.12303626869794793
$ws.Range("J10").Value2 = 0.12382176125828906
$ws.Range("M10").Value2 = 25.0501465
$ws.Range("N10").Value2 = 50.100293
$ws.Range("O10").Value2 = 0.0874724982879541
$ws.Range("P10").Value2 = 0.06006638442832619
$ws.Range("Q10").Value2 = 351.99477216018124
$ws.Range("R10").Value2 = 2111.9686329610877
$ws.Range("S10").Value2 = 0.01076228980303751
$ws.Range("T10").Value2 = 0.007437525512332816

# Row 11
$ws.Range("G11").Value2 = 14.05160533333333
$ws.Range("H11").Value2 = 42.154816
$ws.Range("I11").Value2 = 0.12303626869794793
$ws.Range("J11").Value2 = 0.12382176125828906
$ws.Range("M11").Value2 = 95.05788666666668
$ws.Range("N11").Value2 = 285.17366
$ws.Range("O11").Value2 = 0.3319322235782747
$ws.Range("P11").Value2 = 0.3419012078510756
$ws.Range("Q11").Value2 = 1335.7159072607287
$ws.Range("R11").Value2 = 12021.443165346558
$ws.Range("S11").Value2 = 0.04083970224968393
$ws.Range("T11").Value2 = 0.04233480973245654

# Row 12
$ws.Range("G12").Value2 = 21.412221
$ws.Range("H12").Value2 = 64.236663
$ws.Range("I12").Value2 = 0.1874860355013181
$ws.Range("J12").Value2 = 0.18868299057491247
$ws.Range("M12").Value2 = 13.89934866666667
$ws.Range("N12").Value2 = 41.69804600000001
$ws.Range("O12").Value2 = 0.048535075531341794
$ws.Range("P12").Value2 = 0.04999273878390351
$ws.Range("Q12").Value2 = 297.61592540672206
$ws.Range("R12").Value2 = 2678.5433286604984
$ws.Range("S12").Value2 = 0.009099648894128303
$ws.Range("T12").Value2 = 0.009432779460777328

# Row 13
$ws.Range("G13").Value2 = 21.412221
$ws.Range("H13").Value2 = 64.236663
$ws.Range("I13").Value2 = 0.1874860355013181
$ws.Range("J13").Value2 = 0.18868299057491247
$ws.Range("M13").Value2 = 70.36235166666667
$ws.Range("N13").Value2 = 211.087055
$ws.Range("O13").Value2 = 0.24569799165441703
$ws.Range("P13").Value2 = 0.253077086664408
$ws.Range("Q13").Value2 = 1506.614223966385
$ws.Range("R13").Value2 = 13559.528015697462
$ws.Range("S13").Value2 = 0.04606494238592258
$ws.Range("T13").Value2 = 0.0477513415578268

# Row 14
$ws.Range("G14").Value2 = 21.412221
$ws.Range("H14").Value2 = 64.236663
$ws.Range("I14").Value2 = 0.1874860355013181
$ws.Range("J14").Value2 = 0.18868299057491247
$ws.Range("M14").Value2 = 82.007665
$ws.Range("N14").Value2 = 246.022995
$ws.Range("O14").Value2 = 0.28636221094801234
$ws.Range("P14").Value2 = 0.2949625822722868
$ws.Range("Q14").Value2 = 1755.9662466739649
$ws.Range("R14").Value2 = 15803.696220065684
$ws.Range("S14").Value2 = 0.05368891564803498
$ws.Range("T14").Value2 = 0.05565442213083373

# Row 15
$ws.Range("G15").Value2 = 21.412221
$ws.Range("H15").Value2 = 64.236663
$ws.Range("I15").Value2 = 0.1874860355013181
$ws.Range("J15").Value2 = 0.18868299057491247
$ws.Range("M15").Value2 = 25.0501465
$ws.Range("N15").Value2 = 50.100293
$ws.Range("O15").Value2 = 0.0874724982879541
$ws.Range("P15").Value2 = 0.06006638442832619
$ws.Range("Q15").Value2 = 536.3792729403765
$ws.Range("R15").Value2 = 3218.2756376422585
$ws.Range("S15").Value2 = 0.016399871919404348
$ws.Range("T15").Value2 = 0.01133350504695894

# Row 16
$ws.Range("G16").Value2 = 21.412221
$ws.Range("H16").Value2 = 64.236663
$ws.Range("I16").Value2 = 0.1874860355013181
$ws.Range("J16").Value2 = 0.18868299057491247
$ws.Range("M16").Value2 = 95.05788666666668
$ws.Range("N16").Value2 = 285.17366
$ws.Range("O16").Value2 = 0.3319322235782747
$ws.Range("P16").Value2 = 0.3419012078510756
$ws.Range("Q16").Value2 = 2035.40047709962
$ws.Range("R16").Value2 = 18318.604293896577
$ws.Range("S16").Value2 = 0.06223265665382786
$ws.Range("T16").Value2 = 0.06451094237851568

# Row 17
$ws.Range("G17").Value2 = 2.1734975
$ws.Range("H17").Value2 = 4.346995
$ws.Range("I17").Value2 = 0.019031207899779575
$ws.Range("J17").Value2 = 0.012768471746643995
$ws.Range("M17").Value2 = 13.89934866666667
$ws.Range("N17").Value2 = 41.69804600000001
$ws.Range("O17").Value2 = 0.048535075531341794
$ws.Range("P17").Value2 = 0.04999273878390351
$ws.Range("Q17").Value2 = 30.210199578628337
$ws.Range("R17").Value2 = 181.26119747177003
$ws.Range("S17").Value2 = 0.0009236811128684702
$ws.Range("T17").Value2 = 0.0006383308726996256

# Row 18
$ws.Range("G18").Value2 = 2.1734975
$ws.Range("H18").Value2 = 4.346995
$ws.Range("I18").Value2 = 0.019031207899779575
$ws.Range("J18").Value2 = 0.012768471746643995
$ws.Range("M18").Value2 = 70.36235166666667
$ws.Range("N18").Value2 = 211.087055
$ws.Range("O18").Value2 = 0.24569799165441703
$ws.Range("P18").Value2 = 0.253077086664408
$ws.Range("Q18").Value2 = 152.93239544162083
$ws.Range("R18").Value2 = 917.594372649725
$ws.Range("S18").Value2 = 0.004675929559733517
$ws.Range("T18").Value2 = 0.0032314076307974676

# Row 19
$ws.Range("G19").Value2 = 2.1734975
$ws.Range("H19").Value2 = 4.346995
$ws.Range("I19").Value2 = 0.019031207899779575
$ws.Range("J19").Value2 = 0.012768471746643995
$ws.Range("M19").Value2 = 82.007665
$ws.Range("N19").Value2 = 246.022995
$ws.Range("O19").Value2 = 0.28636221094801234
$ws.Range("P19").Value2 = 0.2949625822722868
$ws.Range("Q19").Value2 = 178.2434548583375
$ws.Range("R19").Value2 = 1069.460729150025
$ws.Range("S19").Value2 = 0.0054498187711921565
$ws.Range("T19").Value2 = 0.0037662213980608483

# Row 20
$ws.Range("G20").Value2 = 2.1734975
$ws.Range("H20").Value2 = 4.346995
$ws.Range("I20").Value2 = 0.019031207899779575
$ws.Range("J20").Value2 = 0.012768471746643995
$ws.Range("M20").Value2 = 25.0501465
$ws.Range("N20").Value2 = 50.100293
$ws.Range("O20").Value2 = 0.0874724982879541
$ws.Range("P20").Value2 = 0.06006638442832619
$ws.Range("Q20").Value2 = 54.446430792383744
$ws.Range("R20").Value2 = 217.78572316953498
$ws.Range("S20").Value2 = 0.001664707300431167
$ws.Range("T20").Value2 = 0.0007669559324961397

# Row 21
$ws.Range("G21").Value2 = 2.1734975
$ws.Range("H21").Value2 = 4.346995
$ws.Range("I21").Value2 = 0.019031207899779575
$ws.Range("J21").Value2 = 0.012768471746643995
$ws.Range("M21").Value2 = 95.05788666666668
$ws.Range("N21").Value2 = 285.17366
$ws.Range("O21").Value2 = 0.3319322235782747
$ws.Range("P21").Value2 = 0.3419012078510756
$ws.Range("Q21").Value2 = 206.60807902528333
$ws.Range("R21").Value2 = 1239.6484741517
$ws.Range("S21").Value2 = 0.00631707115555426
$ws.Range("T21").Value2 = 0.004365555912589915

# Row 22
$ws.Range("G22").Value2 = 28.477822
$ws.Range("H22").Value2 = 85.433466
$ws.Range("I22").Value2 = 0.2493526452249964
$ws.Range("J22").Value2 = 0.2509445713277495
$ws.Range("M22").Value2 = 13.89934866666667
$ws.Range("N22").Value2 = 41.69804600000001
$ws.Range("O22").Value2 = 0.048535075531341794
$ws.Range("P22").Value2 = 0.04999273878390351
$ws.Range("Q22").Value2 = 395.82317724527076
$ws.Range("R22").Value2 = 3562.408595207437
$ws.Range("S22").Value2 = 0.012102349469935074
$ws.Range("T22").Value2 = 0.012545406403626825

# Row 23
$ws.Range("G23").Value2 = 28.477822
$ws.Range("H23").Value2 = 85.433466
$ws.Range("I23").Value2 = 0.2493526452249964
$ws.Range("J23").Value2 = 0.2509445713277495
$ws.Range("M23").Value2 = 70.36235166666667
$ws.Range("N23").Value2 = 211.087055
$ws.Range("O23").Value2 = 0.24569799165441703
$ws.Range("P23").Value2 = 0.253077086664408
$ws.Range("Q23").Value2 = 2003.7665262647367
$ws.Range("R23").Value2 = 18033.898736382627
$ws.Range("S23").Value2 = 0.06126544414549797
$ws.Range("T23").Value2 = 0.06350832102587557

# Row 24
$ws.Range("G24").Value2 = 28.477822
$ws.Range("H24").Value2 = 85.433466
$ws.Range("I24").Value2 = 0.2493526452249964
$ws.Range("J24").Value2 = 0.2509445713277495
$ws.Range("M24").Value2 = 82.007665
$ws.Range("N24").Value2 = 246.022995
$ws.Range("O24").Value2 = 0.28636221094801234
$ws.Range("P24").Value2 = 0.2949625822722868
$ws.Range("Q24").Value2 = 2335.39968650563
$ws.Range("R24").Value2 = 21018.59717855067
$ws.Range("S24").Value2 = 0.0714051747923653
$ws.Range("T24").Value2 = 0.07401925876604504

# Row 25
$ws.Range("G25").Value2 = 28.477822
$ws.Range("H25").Value2 = 85.433466
$ws.Range("I25").Value2 = 0.2493526452249964
$ws.Range("J25").Value2 = 0.2509445713277495
$ws.Range("M25").Value2 = 25.0501465
$ws.Range("N25").Value2 = 50.100293
$ws.Range("O25").Value2 = 0.0874724982879541
$ws.Range("P25").Value2 = 0.06006638442832619
$ws.Range("Q25").Value2 = 713.373613100923
$ws.Range("R25").Value2 = 4280.241678605538
$ws.Range("S25").Value2 = 0.02181149883254032
$ws.Range("T25").Value2 = 0.015073333091574125

# Row 26
$ws.Range("G26").Value2 = 28.477822
$ws.Range("H26").Value2 = 85.433466
$ws.Range("I26").Value2 = 0.2493526452249964
$ws.Range("J26").Value2 = 0.2509445713277495
$ws.Range("M26").Value2 = 95.05788666666668
$ws.Range("N26").Value2 = 285.17366
$ws.Range("O26").Value2 = 0.3319322235782747
$ws.Range("P26").Value2 = 0.3419012078510756
$ws.Range("Q26").Value2 = 2707.041576189507
$ws.Range("R26").Value2 = 24363.374185705557
$ws.Range("S26").Value2 = 0.0827681779846577
$ws.Range("T26").Value2 = 0.08579825204062794
